$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# The 4 existing mailto hyperlinks that live in the block of rows that is
# about to be pushed down by the insert below. Excel's row-insert does not
# repoint pre-existing Hyperlink.Range anchors, so remember their target
# addresses (keyed by their *current* row) up front and fix them up by hand
# afterwards.
$shiftedLinks = @{
    30 = "mailto:auannotationuser2@mailinator.com"
    31 = "mailto:auannotationuser3@mailinator.com"
    32 = "mailto:shareannotationuser1@mailinator.com"
    33 = "mailto:myShareAnnotationUser@mailinator.com"
}

# Insert two new blank rows before the existing row 28 (the blank separator
# row that precedes the "auAnnotationUser*" / "shareAnnotationUser*" block),
# shifting that whole block down to rows 30-36.
$ws.Range("28:29").Insert()

# Repoint the 4 hyperlinks that used to sit at rows 30-33 (now 32-35): drop
# the stale entry (still anchored to the pre-shift row) and re-add it two
# rows further down, pointing at the same mailbox as before.
foreach ($oldRow in $shiftedLinks.Keys) {
    $newRow = $oldRow + 2
    $target = $shiftedLinks[$oldRow]
    $ws.Range("G$oldRow").Hyperlinks.Item(1).Delete()
    $ws.Hyperlinks.Add($ws.Range("G$newRow"), $target)
}

# Row 28: AUtestuser11
$ws.Cells.Item(28, 1).Value = "AUtestuser11"
$ws.Cells.Item(28, 2).Value = "Password1"
$ws.Cells.Item(28, 7).Value = "AUtestuser11@mailinator.com"
$ws.Cells.Item(28, 8).Value = "thomsonreuters"

# Row 29: AUtestuser12
$ws.Cells.Item(29, 1).Value = "AUtestuser12"
$ws.Cells.Item(29, 2).Value = "Password1"
$ws.Cells.Item(29, 7).Value = "AUtestuser12@mailinator.com"
$ws.Cells.Item(29, 8).Value = "thomsonreuters"

# Match the formatting used by the row above (row 27, the previous last
# "AUtestuser" row) - thin box borders around A:G.
$ws.Range("A28:G29").Borders.LineStyle = 1
$ws.Range("A28:G29").Borders.Weight = 2

# Hyperlink the new email cells, same as the existing AUtestuser rows.
$ws.Hyperlinks.Add($ws.Range("G28"), "mailto:AUtestuser11@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G29"), "mailto:AUtestuser12@mailinator.com")

# Update the saved selection to match the author's final cursor position.
$ws.Range("G18").Select()
